# Update forecast-error table: recompute values for rows Q0..Q8 (rows 2-10)
# and add a new row for Q9 (row 11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Q0)
$ws.Range("B2").Value = 0.01778118981931784
$ws.Range("C2").Value = 0.8812756614903005
$ws.Range("D2").Value = 2.316443170155194
$ws.Range("E2").Value = 1.521986586719868
$ws.Range("F2").Value = 1.537026199982511
$ws.Range("G2").Value = 51

# Row 3 (Q1)
$ws.Range("B3").Value = 0.2232283752735154
$ws.Range("C3").Value = 0.8890070218121189
$ws.Range("D3").Value = 2.324356938520629
$ws.Range("E3").Value = 1.524584185448816
$ws.Range("F3").Value = 1.5234647779273
$ws.Range("G3").Value = 50

# Row 4 (Q2)
$ws.Range("B4").Value = 0.07282520553819657
$ws.Range("C4").Value = 0.8920146246384424
$ws.Range("D4").Value = 2.369955897058822
$ws.Range("E4").Value = 1.539466107798032
$ws.Range("F4").Value = 1.553678209244919
$ws.Range("G4").Value = 49

# Row 5 (Q3)
$ws.Range("B5").Value = 0.2280708590847791
$ws.Range("C5").Value = 0.9180197143133344
$ws.Range("D5").Value = 2.436291467276783
$ws.Range("E5").Value = 1.5608624113857
$ws.Range("F5").Value = 1.560450066177918
$ws.Range("G5").Value = 48

# Row 6 (Q4)
$ws.Range("B6").Value = 0.1059889376669026
$ws.Range("C6").Value = 0.8775097619739644
$ws.Range("D6").Value = 2.394266195833348
$ws.Range("E6").Value = 1.547341654526675
$ws.Range("F6").Value = 1.560396615019112
$ws.Range("G6").Value = 47

# Row 7 (Q5)
$ws.Range("B7").Value = 0.2424232887606824
$ws.Range("C7").Value = 0.8995544753595199
$ws.Range("D7").Value = 2.464361956921573
$ws.Range("E7").Value = 1.569828639349395
$ws.Range("F7").Value = 1.56813600370666
$ws.Range("G7").Value = 46

# Row 8 (Q6)
$ws.Range("B8").Value = 0.08228634011495191
$ws.Range("C8").Value = 0.8063847497337052
$ws.Range("D8").Value = 2.322743373247686
$ws.Range("E8").Value = 1.524054911493574
$ws.Range("F8").Value = 1.539028286382188
$ws.Range("G8").Value = 45

# Row 9 (Q7)
$ws.Range("B9").Value = 0.2065927448523342
$ws.Range("C9").Value = 0.8408213022771569
$ws.Range("D9").Value = 2.429435650446118
$ws.Range("E9").Value = 1.558664701097102
$ws.Range("F9").Value = 1.562773502624323
$ws.Range("G9").Value = 44

# Row 10 (Q8) - F10 is newly added here
$ws.Range("B10").Value = 0.1150341483697026
$ws.Range("C10").Value = 0.7847972820185366
$ws.Range("D10").Value = 2.399366321466991
$ws.Range("E10").Value = 1.548988806114167
$ws.Range("F10").Value = 1.562992696000677
$ws.Range("G10").Value = 43

# Row 11 (Q9) - new row
$ws.Range("A11").Value = "Q9"
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("B11").Value = 0.124196579601926
$ws.Range("C11").Value = 0.8001464258099134
$ws.Range("D11").Value = 2.371409906460054
$ws.Range("E11").Value = 1.539938280081398
$ws.Range("F11").Value = 1.553527652694611
$ws.Range("G11").Value = 42
